# Updated UI, added logic for handling color click
#
# Typography sheet (sheet 1): the Fallback Character column (F) for rows
# 4-6 is cleared, and the Wildcard Characters column (G) for those rows
# now contains "0123456789abcdef#".
#
# Translation sheet (sheet 2): the SPA column is removed for the header
# row (G3), the "Power" text id's SPA translation is removed (G4) and its
# ENG value changes from "Power<value>" to "<value>" (F4). The whole
# "SingleUseId2" row (row 5) is cleared out.

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: clear Fallback Character, set Wildcard Characters ---
$wsTypography.Range("F4").ClearContents()
$wsTypography.Range("G4").Value = "0123456789abcdef#"

$wsTypography.Range("F5").ClearContents()
$wsTypography.Range("G5").Value = "0123456789abcdef#"

$wsTypography.Range("F6").ClearContents()
$wsTypography.Range("G6").Value = "0123456789abcdef#"

# --- Translation sheet: drop SPA column header ---
$wsTranslation.Range("G3").ClearContents()

# --- Translation sheet: update Power row (row 4) values ---
$wsTranslation.Range("F4").Value = "<value>"
$wsTranslation.Range("G4").ClearContents()

# --- Translation sheet: clear out the SingleUseId2 row entirely ---
$wsTranslation.Range("B5:G5").ClearContents()
